$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "Data Source" (column L) to host "Supplier Name".
# This shifts the existing L (Data Source) -> M and M (Very Custom) -> N.
$ws.Columns("L").Insert()

# Header
$ws.Range("L1").Value = "Supplier Name"

# Supplier Name values per row, grouped by supplier
$ws.Range("L2").Value = "Glassy Glass inc."
$ws.Range("L3").Value = "Glassy Glass inc."
$ws.Range("L4").Value = "Spice girls inc."
$ws.Range("L5").Value = "Spice girls inc."
$ws.Range("L6").Value = "Spice girls inc."
$ws.Range("L7").Value = "Spice girls inc."
$ws.Range("L8").Value = "Ship happens inc."
$ws.Range("L9").Value = "Ship happens inc."
$ws.Range("L10").Value = "Dumpster Divers inc."
$ws.Range("L11").Value = "Dumpster Divers inc."
